# Scheduled-runner update: refresh cached market-price / leve-profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the ALC, ARM, BSM,
# CRP, CUL and GSM sheets. Values below are the new cached numbers; a few
# rows also gain/lose an N (LeveProfitHQ) or M (LeveProfitNQ) cell because
# that column had no cached value for that leve before/after this refresh.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H32").Value = 1600
$ws_ALC.Range("I32").Value = 1250.25
$ws_ALC.Range("J32").Value = 1799.8572
$ws_ALC.Range("K32").Value = 1250.25
$ws_ALC.Range("L32").Value = 1799.8572
$ws_ALC.Range("M32").Value = -924.25
$ws_ALC.Range("N32").Value = -2451.8572
$ws_ALC.Range("H33").Value = 808.3333
$ws_ALC.Range("I33").Value = 1024.6875
$ws_ALC.Range("J33").Value = 116
$ws_ALC.Range("K33").Value = 1024.6875
$ws_ALC.Range("L33").Value = 116
$ws_ALC.Range("M33").Value = -795.6875
$ws_ALC.Range("N33").Value = -574
$ws_ALC.Range("H40").Value = 940
$ws_ALC.Range("I40").Value = 975
$ws_ALC.Range("J40").Value = 800
$ws_ALC.Range("K40").Value = 975
$ws_ALC.Range("L40").Value = 800
$ws_ALC.Range("M40").Value = -800
$ws_ALC.Range("N40").Value = -1150
$ws_ALC.Range("H51").Value = 31670.1
$ws_ALC.Range("I51").Value = 50833
$ws_ALC.Range("J51").Value = 23457.428
$ws_ALC.Range("K51").Value = 50833
$ws_ALC.Range("L51").Value = 23457.428
$ws_ALC.Range("M51").Value = -50349
$ws_ALC.Range("N51").Value = -24425.428
$ws_ALC.Range("H55").Value = 312.91666
$ws_ALC.Range("I55").Value = 100.625
$ws_ALC.Range("J55").Value = 737.5
$ws_ALC.Range("K55").Value = 100.625
$ws_ALC.Range("L55").Value = 737.5
$ws_ALC.Range("M55").Value = 113.375
$ws_ALC.Range("N55").Value = -1165.5
$ws_ALC.Range("H64").Value = 3506.25
$ws_ALC.Range("I64").Value = 3225.2632
$ws_ALC.Range("J64").Value = 4574
$ws_ALC.Range("K64").Value = 3225.2632
$ws_ALC.Range("L64").Value = 4574
$ws_ALC.Range("M64").Value = -2977.2632
$ws_ALC.Range("N64").Value = -5070
$ws_ALC.Range("H67").Value = 3506.25
$ws_ALC.Range("I67").Value = 3225.2632
$ws_ALC.Range("J67").Value = 4574
$ws_ALC.Range("K67").Value = 3225.2632
$ws_ALC.Range("L67").Value = 4574
$ws_ALC.Range("M67").Value = -2367.2632
$ws_ALC.Range("N67").Value = -6290
$ws_ALC.Range("H76").Value = 333340930
$ws_ALC.Range("I76").Value = 500010000
$ws_ALC.Range("J76").Value = 2800
$ws_ALC.Range("K76").Value = 500010000
$ws_ALC.Range("L76").Value = 2800
$ws_ALC.Range("M76").Value = -500009685
$ws_ALC.Range("N76").Value = -3430
$ws_ALC.Range("H79").Value = 333340930
$ws_ALC.Range("I79").Value = 500010000
$ws_ALC.Range("J79").Value = 2800
$ws_ALC.Range("K79").Value = 500010000
$ws_ALC.Range("L79").Value = 2800
$ws_ALC.Range("M79").Value = -500008908
$ws_ALC.Range("N79").Value = -4984
$ws_ALC.Range("H80").Value = 3382080
$ws_ALC.Range("I80").Value = 6480.421
$ws_ALC.Range("J80").Value = 6945213
$ws_ALC.Range("K80").Value = 19441.263
$ws_ALC.Range("L80").Value = 20835639
$ws_ALC.Range("M80").Value = -18443.263
$ws_ALC.Range("N80").Value = -20837635
$ws_ALC.Range("H83").Value = 3382080
$ws_ALC.Range("I83").Value = 6480.421
$ws_ALC.Range("J83").Value = 6945213
$ws_ALC.Range("K83").Value = 58323.789
$ws_ALC.Range("L83").Value = 62506917
$ws_ALC.Range("M83").Value = -53331.789
$ws_ALC.Range("N83").Value = -62516901
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H4").Value = 92.85714
$ws_ARM.Range("I4").Value = 93.333336
$ws_ARM.Range("J4").Value = 90
$ws_ARM.Range("K4").Value = 93.333336
$ws_ARM.Range("L4").Value = 90
$ws_ARM.Range("M4").Value = 22.666664
$ws_ARM.Range("N4").Value = -322
$ws_ARM.Range("H63").Value = 2249.875
$ws_ARM.Range("I63").Value = 2266.6667
$ws_ARM.Range("J63").Value = 2239.8
$ws_ARM.Range("K63").Value = 2266.6667
$ws_ARM.Range("L63").Value = 2239.8
$ws_ARM.Range("M63").Value = -1580.6667
$ws_ARM.Range("N63").Value = -3611.8
$ws_ARM.Range("H66").Value = 2249.875
$ws_ARM.Range("I66").Value = 2266.6667
$ws_ARM.Range("J66").Value = 2239.8
$ws_ARM.Range("K66").Value = 11333.3335
$ws_ARM.Range("L66").Value = 11199
$ws_ARM.Range("M66").Value = -7901.333500000001
$ws_ARM.Range("N66").Value = -18063
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H26").Value = 7133.3335
$ws_BSM.Range("I26").Value = 7133.3335
$ws_BSM.Range("J26").Value = 0
$ws_BSM.Range("K26").Value = 7133.3335
$ws_BSM.Range("L26").Value = 0
$ws_BSM.Range("M26").Value = -6841.3335
$ws_BSM.Range("N26").ClearContents()
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H22").Value = 813.5294
$ws_CRP.Range("I22").Value = 194.16667
$ws_CRP.Range("J22").Value = 2300
$ws_CRP.Range("K22").Value = 194.16667
$ws_CRP.Range("L22").Value = 2300
$ws_CRP.Range("M22").Value = 155.83333
$ws_CRP.Range("N22").Value = -3000
$ws_CRP.Range("H62").Value = 3101.1
$ws_CRP.Range("I62").Value = 2226.25
$ws_CRP.Range("J62").Value = 3684.3333
$ws_CRP.Range("K62").Value = 2226.25
$ws_CRP.Range("L62").Value = 3684.3333
$ws_CRP.Range("M62").Value = -1602.25
$ws_CRP.Range("N62").Value = -4932.3333
$ws_CRP.Range("H65").Value = 3101.1
$ws_CRP.Range("I65").Value = 2226.25
$ws_CRP.Range("J65").Value = 3684.3333
$ws_CRP.Range("K65").Value = 11131.25
$ws_CRP.Range("L65").Value = 18421.6665
$ws_CRP.Range("M65").Value = -8011.25
$ws_CRP.Range("N65").Value = -24661.6665
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 2084848.9
$ws_CUL.Range("I5").Value = 692.25
$ws_CUL.Range("J5").Value = 3474286.8
$ws_CUL.Range("K5").Value = 2076.75
$ws_CUL.Range("L5").Value = 10422860.4
$ws_CUL.Range("M5").Value = -1964.75
$ws_CUL.Range("N5").Value = -10423084.4
$ws_CUL.Range("H23").Value = 10087.5
$ws_CUL.Range("I23").Value = 82.5
$ws_CUL.Range("J23").Value = 16757.5
$ws_CUL.Range("K23").Value = 247.5
$ws_CUL.Range("L23").Value = 50272.5
$ws_CUL.Range("M23").Value = -12.5
$ws_CUL.Range("N23").Value = -50742.5
$ws_CUL.Range("H39").Value = 3712.25
$ws_CUL.Range("I39").Value = 700.5
$ws_CUL.Range("J39").Value = 4142.5
$ws_CUL.Range("K39").Value = 2101.5
$ws_CUL.Range("L39").Value = 12427.5
$ws_CUL.Range("M39").Value = -1807.5
$ws_CUL.Range("N39").Value = -13015.5
$ws_CUL.Range("H86").Value = 493.33334
$ws_CUL.Range("I86").Value = 0
$ws_CUL.Range("J86").Value = 493.33334
$ws_CUL.Range("K86").Value = 0
$ws_CUL.Range("L86").Value = 1480.00002
$ws_CUL.Range("M86").ClearContents()
$ws_CUL.Range("N86").Value = -3852.00002
$ws_CUL.Range("H89").Value = 493.33334
$ws_CUL.Range("I89").Value = 0
$ws_CUL.Range("J89").Value = 493.33334
$ws_CUL.Range("K89").Value = 0
$ws_CUL.Range("L89").Value = 4440.00006
$ws_CUL.Range("M89").ClearContents()
$ws_CUL.Range("N89").Value = -16296.00006
$ws_CUL.Range("H131").Value = 880.7
$ws_CUL.Range("I131").Value = 360.69232
$ws_CUL.Range("J131").Value = 958.4023
$ws_CUL.Range("K131").Value = 1082.07696
$ws_CUL.Range("L131").Value = 2875.2069
$ws_CUL.Range("M131").Value = 3957.92304
$ws_CUL.Range("N131").Value = -12955.2069
$ws_CUL.Range("H135").Value = 2084848.9
$ws_CUL.Range("I135").Value = 692.25
$ws_CUL.Range("J135").Value = 3474286.8
$ws_CUL.Range("K135").Value = 6230.25
$ws_CUL.Range("L135").Value = 31268581.2
$ws_CUL.Range("M135").Value = -3695.25
$ws_CUL.Range("N135").Value = -31273651.2
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H70").Value = 14672824
$ws_GSM.Range("I70").Value = 5213116.5
$ws_GSM.Range("J70").Value = 71431070
$ws_GSM.Range("K70").Value = 5213116.5
$ws_GSM.Range("L70").Value = 71431070
$ws_GSM.Range("M70").Value = -5212846.5
$ws_GSM.Range("N70").Value = -71431610
$ws_GSM.Range("H73").Value = 14672824
$ws_GSM.Range("I73").Value = 5213116.5
$ws_GSM.Range("J73").Value = 71431070
$ws_GSM.Range("K73").Value = 5213116.5
$ws_GSM.Range("L73").Value = 71431070
$ws_GSM.Range("M73").Value = -5212180.5
$ws_GSM.Range("N73").Value = -71432942
$ws_GSM.Range("H80").Value = 7135.122
$ws_GSM.Range("I80").Value = 3859.3333
$ws_GSM.Range("J80").Value = 16069.091
$ws_GSM.Range("K80").Value = 3859.3333
$ws_GSM.Range("L80").Value = 16069.091
$ws_GSM.Range("M80").Value = -2861.3333
$ws_GSM.Range("N80").Value = -18065.091
$ws_GSM.Range("H83").Value = 7135.122
$ws_GSM.Range("I83").Value = 3859.3333
$ws_GSM.Range("J83").Value = 16069.091
$ws_GSM.Range("K83").Value = 19296.6665
$ws_GSM.Range("L83").Value = 80345.455
$ws_GSM.Range("M83").Value = -14304.6665
$ws_GSM.Range("N83").Value = -90329.455
$ws_GSM.Range("H136").Value = 22600.295
$ws_GSM.Range("I136").Value = 0
$ws_GSM.Range("J136").Value = 22600.295
$ws_GSM.Range("K136").Value = 0
$ws_GSM.Range("L136").Value = 67800.885
$ws_GSM.Range("N136").Value = -72900.885
